$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert a new row 4 in the "fragments" sheet for the new "auth.email" fragment
$ws1.Rows.Item(4).Insert()
$ws1.Range("A4").Value = "auth.email"
$ws1.Range("B4").Value = 0
$ws1.Range("D4").Value = "E-mail"
$ws1.Range("E4").Value = "E-mail"

# Extend the _FilterDatabase defined name to include the new row
$wb.Names.Item("fragments!_FilterDatabase").RefersTo = "=fragments!`$A`$2:`$Z`$46"

# Update view/selection state: select a cell on the "hidden" sheet first ...
$ws2.Select()
$ws2.Range("B98").Select()

# ... then make "fragments" the active sheet with A6 selected
$ws1.Select()
$ws1.Range("A6").Select()
